# "segunda parte de los articulos"
# Remove the leftover "IA:" / "IA Coversacional : GPT-4.5" scratch notes
# (B9:B10) that were left over from the first part of the assignment, now
# that the second half of the article list has been filled in below.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").ClearContents() | Out-Null
$ws.Range("B10").ClearContents() | Out-Null

# Leave the cursor where the author ended up after cleaning the sheet.
$ws.Range("E13").Select() | Out-Null
